# Update LDLC prices history
# Insert a new price-snapshot column right before the "nom" column (GD),
# shifting "nom" (GD->GE) and "url_produit" (GE->GF) one column to the right.
# The new column gets the latest timestamp as its header (row 1) and, for
# each product row, a copy of the previous last snapshot value (old GC,
# now GC still) so the new column mirrors the most recent known price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GD is column 186 (1-based). Insert a new column there; everything from
# GD onward (nom, url_produit) shifts right by one column.
$ws.Range("GD1").EntireColumn.Insert()

# New column header: latest snapshot timestamp.
$ws.Range("GD1").Value2 = "2026-02-05 14:40:21"

# Populate the new GD column (rows 2-210) with the same value as the GC
# column on that row (the previous most-recent snapshot). Rows where GC
# is blank are left untouched (already blank after the column insert).
for ($r = 2; $r -le 210; $r++) {
    $gc = $ws.Cells.Item($r, 185)
    $v = $gc.Value2
    if ($null -ne $v -and $v -ne "") {
        $ws.Cells.Item($r, 186).Value2 = $v
    }
}
